# Fruta / hortaliza, semanal
#
# The published update re-sorted/shuffled the 40 data rows (rows 2-41) of
# the sheet. No cell values were actually changed - every row in the
# "after" version is identical to one of the rows from the "before"
# version, just relocated to a different row number. This script captures
# that row permutation and re-writes the data block in the new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source data block (row 2 .. row 41, columns A..T)
$srcRange = $ws.Range("A2:T41")
$srcValues = $srcRange.Value2

# For destination array row i (1-based, i=1 => worksheet row 2), this is the
# 1-based row index into $srcValues that should be copied there.
$rowMap = @(10,4,29,19,24,28,36,11,13,23,30,1,39,16,6,37,15,3,25,38,31,26,18,14,21,2,12,35,9,20,8,34,22,5,33,7,32,17,40,27)

$rowCount = $srcValues.GetLength(0)
$colCount = $srcValues.GetLength(1)

$newValues = New-Object 'object[,]' $rowCount, $colCount

for ($i = 1; $i -le $rowCount; $i++) {
    $srcRow = $rowMap[$i - 1]
    for ($c = 1; $c -le $colCount; $c++) {
        $newValues[$i - 1, $c - 1] = $srcValues[$srcRow, $c]
    }
}

$srcRange.Value2 = $newValues
